$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be
# explicitly formatted as Text first, otherwise Excel auto-converts
# the inline-string price into a numeric value.
$textForceCells = @(
    'D4',
    'D5',
    'D7',
    'D9',
    'D10',
    'D11',
    'D12',
    'D14',
    'D15',
    'D16',
    'D19',
    'D21',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D48',
    'D49',
    'D50',
    'D51',
)

foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$updates = @{
    'D2' = '30.570.86'
    'E2' = '  -1.51%  '
    'D3' = '1.922.18'
    'E3' = '  -1.91%  '
    'D4' = '1.003'
    'E4' = '  +0.37%  '
    'D5' = '239.24'
    'E5' = '  -2.68%  '
    'E6' = '  +0.33%  '
    'D7' = '0.4801'
    'E7' = '  -1.71%  '
    'E8' = '  -2.94%  '
    'D9' = '0.06710'
    'E9' = '  -1.84%  '
    'D10' = '18.79'
    'E10' = '  -2.30%  '
    'D11' = '103.90'
    'E11' = '  -2.90%  '
    'B12' = 'TRON'
    'C12' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D12' = '0.07747'
    'E12' = '  -1.15%  '
    'B13' = 'WrappedEther'
    'C13' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D13' = '1.919.52'
    'E13' = '  -2.21%  '
    'D14' = '5.235'
    'E14' = '  -4.93%  '
    'D15' = '0.6818'
    'D16' = '265.81'
    'E16' = '  -6.58%  '
    'D17' = '30.613.19'
    'E17' = '  -1.45%  '
    'D19' = '0.000007539'
    'E19' = '  -2.29%  '
    'D21' = '5.432'
    'E21' = '  -1.80%  '
    'E22' = '  +0.53%  '
    'B23' = 'Chainlink'
    'C23' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D23' = '6.342'
    'E23' = '  -2.79%  '
    'B24' = 'Cosmos'
    'C24' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D24' = '9.633'
    'E24' = '  -2.06%  '
    'B25' = 'Monero'
    'C25' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D25' = '163.37'
    'E25' = '  -3.39%  '
    'B26' = 'EthereumClassic'
    'C26' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D26' = '19.06'
    'E26' = '  -4.91%  '
    'B27' = 'LidoDAOToken'
    'C27' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D27' = '2.100'
    'E27' = '  -4.97%  '
    'B28' = 'Stellar'
    'C28' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D28' = '0.1024'
    'E28' = '  -3.22%  '
    'B29' = 'Toncoin'
    'C29' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D29' = '1.384'
    'E29' = '  -0.87%  '
    'B30' = 'Filecoin'
    'C30' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D30' = '4.570'
    'E30' = '  -0.90%  '
    'B31' = 'PancakeSwap'
    'C31' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D31' = '1.517'
    'E31' = '  -4.39%  '
    'B32' = 'InternetComputer(DFINITY)'
    'C32' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D32' = '4.265'
    'E32' = '  -4.21%  '
    'B33' = 'Hedera'
    'C33' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D33' = '0.04757'
    'E33' = '  -3.93%  '
    'B34' = 'ImmutableX'
    'C34' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D34' = '0.7382'
    'E34' = '  -3.41%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D35' = '1.124'
    'E35' = '  -4.44%  '
    'B36' = 'Frax'
    'C36' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'D36' = '1.003'
    'E36' = '  +0.29%  '
    'B37' = 'HuobiToken'
    'C37' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D37' = '2.684'
    'E37' = '  -1.78%  '
    'B38' = 'VeChain'
    'C38' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D38' = '0.01942'
    'E38' = '  -4.25%  '
    'B39' = 'MXToken'
    'C39' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D39' = '2.640'
    'E39' = '  -2.29%  '
    'B40' = 'FraxShare'
    'C40' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D40' = '6.353'
    'E40' = '  -3.25%  '
    'B41' = 'Aave'
    'C41' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D41' = '75.37'
    'E41' = '  -3.85%  '
    'B42' = 'RenderToken'
    'C42' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D42' = '2.009'
    'E42' = '  -5.26%  '
    'B43' = 'TrustWalletToken'
    'C43' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D43' = '0.8617'
    'E43' = '  -5.19%  '
    'B44' = 'Quant'
    'C44' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D44' = '106.45'
    'E44' = '  -2.51%  '
    'B45' = 'TheSandbox'
    'C45' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D45' = '0.4292'
    'E45' = '  -4.41%  '
    'B46' = 'PaxDollar'
    'C46' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D46' = '1.003'
    'E46' = '  +0.24%  '
    'B47' = 'Maker'
    'C47' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D47' = '1.002.08'
    'E47' = '  -2.03%  '
    'D48' = '7.525'
    'E48' = '  -8.05%  '
    'B49' = 'Algorand'
    'C49' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D49' = '0.1204'
    'E49' = '  -4.79%  '
    'D50' = '35.25'
    'E50' = '  -2.12%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D51' = '8.978'
    'E51' = '  -4.63%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
